$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Names"

# Replace Vala/Valic row with Rada/Radic
$ws.Range("A3").Value = "Rada"
$ws.Range("B3").Value = "Radic"

# Update selection to B4
$ws.Range("B4").Select()
